$d = $word.ActiveDocument

# 1) Remove bold formatting from the "Präsentation, Klassendiagramm, ..." run
#    (spans across three runs: "..Setup für Abgabe (lokale Version - " / "web" / " Version - Testläufe), Ablaufplan")
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "Präsentation, Klassendiagramm, Setup für Abgabe (lokale Version - web Version - Testläufe), Ablaufplan"
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
$found = $rng.Find.Execute()
if ($found) {
    $rng.Font.Bold = 0
}

# 2) Remove the lastRenderedPageBreak field before "Zukünftige / "
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Text = "Zukünftige / "
$rng2.Find.Forward = $true
$rng2.Find.Wrap = 0
$found2 = $rng2.Find.Execute()
if ($found2) {
    $collapsed = $rng2.Duplicate
    $collapsed.Collapse(1)
    $collapsed.MoveStart(1, -1)
    $collapsed.Delete()
}
